$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates for rows 2-44 ---
$ws.Range("D2").Value = "'331.99"
$ws.Range("E2").Value = "'1.62%"
$ws.Range("D3").Value = "'45.85"
$ws.Range("E3").Value = "'3.96%"
$ws.Range("D4").Value = "'5.695"
$ws.Range("E4").Value = "'3.15%"
$ws.Range("D5").Value = "'0.08367"
$ws.Range("E5").Value = "'4.26%"
$ws.Range("D6").Value = "'2.041"
$ws.Range("E6").Value = "'2.42%"
$ws.Range("D7").Value = "'4.493"
$ws.Range("E7").Value = "'4.62%"
$ws.Range("D8").Value = "'0.9835"
$ws.Range("E8").Value = "'3.25%"
$ws.Range("D9").Value = "'2.588"
$ws.Range("E9").Value = "'0.54%"
$ws.Range("D10").Value = "'0.1164"
$ws.Range("E10").Value = "'2.80%"
$ws.Range("D11").Value = "'0.1942"
$ws.Range("E11").Value = "'3.98%"
$ws.Range("D12").Value = "'10.40"
$ws.Range("E12").Value = "'-2.13%"
$ws.Range("D13").Value = "'0.1006"
$ws.Range("E13").Value = "'1.35%"
$ws.Range("D14").Value = "'0.04649"
$ws.Range("E14").Value = "'0.61%"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("E15").Value = "'-0.51%"
$ws.Range("D16").Value = "'0.001287"
$ws.Range("E16").Value = "'1.55%"
$ws.Range("D17").Value = "'0.006095"
$ws.Range("E17").Value = "'4.46%"
$ws.Range("D18").Value = "'3.369"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("D19").Value = "'0.3365"
$ws.Range("E19").Value = "'-3.24%"
$ws.Range("D20").Value = "'0.1400"
$ws.Range("E20").Value = "'-0.66%"
$ws.Range("E21").Value = "'2.34%"
$ws.Range("D22").Value = "'0.04210"
$ws.Range("E22").Value = "'3.16%"
$ws.Range("D23").Value = "'0.001308"
$ws.Range("E23").Value = "'4.95%"
$ws.Range("D24").Value = "'0.004685"
$ws.Range("E24").Value = "'8.48%"
$ws.Range("E25").Value = "'7.55%"
$ws.Range("D26").Value = "'0.0003742"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("D38").Value = "'0.02795"
$ws.Range("E38").Value = "'9.26%"
$ws.Range("D39").Value = "'0.05816"
$ws.Range("E39").Value = "'4.27%"
$ws.Range("D40").Value = "'0.007745"
$ws.Range("E40").Value = "'2.39%"
$ws.Range("E41").Value = "'3.14%"
$ws.Range("D42").Value = "'0.007201"
$ws.Range("E42").Value = "'-5.54%"
$ws.Range("E43").Value = "'-2.07%"
$ws.Range("D44").Value = "'0.008102"
$ws.Range("E44").Value = "'-4.85%"
# --- Rows 45-51: a new coin (PooCoin) was inserted at rank 45, shifting
#     CoinLion..CryptobidCoin down by one row; SpecialPowerGold (old row 51) drops off ---
$ws.Range("B45").Value = "PooCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ucHyn6T7W+poocoin-poocoin"
$ws.Range("D45").Value = "'0.3501"
$ws.Range("E45").Value = "'--%"
$ws.Range("B46").Value = "CoinLion"
$ws.Range("C46").Value = "https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion"
$ws.Range("D46").Value = "'0.00007300"
$ws.Range("E46").Value = "'2.55%"
$ws.Range("B47").Value = "Kangarootoken"
$ws.Range("C47").Value = "https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("B48").Value = "ACDXExchange"
$ws.Range("C48").Value = "https://coinranking.com/coin/-y35lbZ7U+acdxexchange-acxt"
$ws.Range("D48").Value = "'0.0005804"
$ws.Range("E48").Value = "'-0.13%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003500"
$ws.Range("E49").Value = "'17.44%"
$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D50").Value = "'0.003500"
$ws.Range("E50").Value = "'-0.81%"
$ws.Range("B51").Value = "CryptobidCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D51").Value = "'0.00002101"
$ws.Range("E51").Value = "'0.13%"
